# "Scrum Master: Files Update"
# - Adds a new backlog card "Pesquisar como se joga o jogo" to the
#   "To do:" column (cell E6), as a new shared string.
# - Scrolls the sheet view so column C is the left-most visible column
#   and selects the newly-added cell E6 (was D4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New backlog item text -> E6 (creates a new shared-string entry).
$ws.Range("E6").Value = "Pesquisar como se joga o jogo"

# Scroll so column C is the top-left visible column, then select E6
# and make it the active cell (matches sheetView/selection in the sheet).
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E6").Select()
